$wb = $excel.ActiveWorkbook

# --- 1. Rename headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet right after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting used on the other sheets: bold font, thin box
# border, centered horizontally, top-aligned vertically.
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45417.99999999999, 54, 32.9416894273066, 75.01771314240953),
    @(45494.99999999999, 37, 16.83170934424722, 56.5838936067297),
    @(45501.99999999999, 35, 15.55758863498953, 56.20110523998164),
    @(45522.99999999999, 31, 9.577811638751943, 50.70835425084447),
    @(45557.99999999999, 23, 3.161244794593854, 42.92045647745899),
    @(45564.99999999999, 21, 0.5795215325971143, 41.49264972489001),
    @(45571.99999999999, 19, -0.3500463161652522, 40.19673295917319),
    @(45592.99999999999, 15, -6.197014232348746, 35.72812345130689),
    @(45599.99999999999, 13, -6.689074590691943, 33.97820986966489),
    @(45606.99999999999, 11, -8.901648780787891, 32.11864437618321),
    @(45613.99999999999, 10, -9.999452412272255, 29.26272255091402),
    @(45620.99999999999, 8, -12.53796453734466, 29.93101369429449),
    @(45627.99999999999, 7, -15.91715551205265, 27.82739519824219),
    @(45634.99999999999, 5, -15.45640109995259, 24.57479265297891),
    @(45641.99999999999, 3, -18.85982336576749, 24.19558575416574),
    @(45648.99999999999, 2, -18.19844612632044, 23.66444811951336),
    @(45655.99999999999, 0, -21.79825858500599, 19.60030796459364)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Column A holds dates - match the date number format used for the
# "Order Week" / "Order Month" columns on the other two sheets.
$wsForecast.Range("A2:A18").NumberFormat = $wsWeekly.Range("A2").NumberFormat
